$d = $word.ActiveDocument

$replacements = @(
    @{old = "57×33=1881"; new = "83×47=3901"},
    @{old = "29×80=2320"; new = "79×38=3002"},
    @{old = "39×62=2418"; new = "70×77=5390"},
    @{old = "41×25=1025"; new = "76×74=5624"},
    @{old = "58×79=4582"; new = "11×49=539"},
    @{old = "38×48=1824"; new = "23×42=966"},
    @{old = "39×31=1209"; new = "65×19=1235"},
    @{old = "88×13=1144"; new = "27×58=1566"},
    @{old = "40×20=800"; new = "75×56=4200"},
    @{old = "94×14=1316"; new = "92×76=6992"},
    @{old = "11×60=660"; new = "53×92=4876"},
    @{old = "53×55=2915"; new = "75×35=2625"},
    @{old = "55×12=660"; new = "44×26=1144"},
    @{old = "21×48=1008"; new = "35×30=1050"},
    @{old = "68×12=816"; new = "18×74=1332"},
    @{old = "86×53=4558"; new = "88×81=7128"},
    @{old = "59×89=5251"; new = "98×79=7742"},
    @{old = "11×32=352"; new = "23×49=1127"},
    @{old = "96×97=9312"; new = "96×41=3936"},
    @{old = "29×91=2639"; new = "67×22=1474"},
    @{old = "90×59=5310"; new = "91×62=5642"},
    @{old = "50×89=4450"; new = "16×53=848"},
    @{old = "77×88=6776"; new = "62×54=3348"},
    @{old = "50×30=1500"; new = "30×40=1200"},
    @{old = "32×59=1888"; new = "16×18=288"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
